# Revert "adding term 2.0 now utf-8"
# This reverts commit 78ae09ef7158944b5ce8ba326bcf0fbc536c742d, rolling the
# workbook's metadata and concept list back to the 1.1.0 release values.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Metadata": roll Version / Date / Contact back to 1.1.0 values
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value  = "1.1.0"                                 # Version
$ws1.Range("B8").Value  = "2023-07-10T23:08:03+02:00"             # Date
$ws1.Range("B10").Value = "No display for ContactDetail"          # Contact

# -----------------------------------------------------------------
# Sheet "Include from SNOMED CT": restore the 1.1.0 concept list
# (codes 160245001 + 116223007 instead of the single 116224001 concept)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Include from SNOMED CT")

# Update the existing concept code in place (force text, then restore the
# surrounding "General" number formatting so the cell style matches its
# siblings again).
$ws2.Range("A2").Value = "'160245001"
$ws2.Range("A3:B3").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Insert a new row for the second concept code, directly above the blank
# separator row that precedes the "System URI" row.
$ws2.Rows.Item(3).Insert()
$ws2.Range("A4:B4").Copy()
$ws2.Range("A3:B3").PasteSpecial(-4122)

$ws2.Range("A3").Value = "'116223007"
$ws2.Range("A4:B4").Copy()
$ws2.Range("A3:B3").PasteSpecial(-4122)
